$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 51595
$ws.Range("J47").Value = 51595
$ws.Range("L47").Value = 51595
$ws.Range("N47").Value = -53539

$ws.Range("H116").Value = 7399.6294
$ws.Range("I116").Value = 7377
$ws.Range("K116").Value = 7377
$ws.Range("M116").Value = -3935

$ws.Range("H132").Value = 1894.1224
$ws.Range("I132").Value = 1791.25
$ws.Range("K132").Value = 5373.75
$ws.Range("M132").Value = -2843.75

$ws.Range("H136").Value = 349993.25
$ws.Range("J136").Value = 349993.25
$ws.Range("L136").Value = 349993.25
$ws.Range("N136").Value = -360193.25

$ws.Range("H137").Value = 1886.5518
$ws.Range("I137").Value = 1705.8948
$ws.Range("K137").Value = 5117.6844
$ws.Range("M137").Value = -2567.6844

$ws.Range("H138").Value = 2420300.2
$ws.Range("I138").Value = 6341.5
$ws.Range("J138").Value = 4633096
$ws.Range("K138").Value = 19024.5
$ws.Range("L138").Value = 13899288
$ws.Range("M138").Value = -13884.5
$ws.Range("N138").Value = -13909568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6326.7666
$ws.Range("I32").Value = 5022.7144
$ws.Range("J32").Value = 12135.728
$ws.Range("K32").Value = 5022.7144
$ws.Range("L32").Value = 12135.728
$ws.Range("M32").Value = -4735.7144
$ws.Range("N32").Value = -12709.728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3408.25
$ws.Range("I99").Value = 2533.2144
$ws.Range("K99").Value = 2533.2144
$ws.Range("M99").Value = -1035.2144

$ws.Range("H105").Value = 6725.3667
$ws.Range("I105").Value = 7826.524
$ws.Range("J105").Value = 4156
$ws.Range("K105").Value = 7826.524
$ws.Range("L105").Value = 4156
$ws.Range("M105").Value = -6079.524
$ws.Range("N105").Value = -7650

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3083.75
$ws.Range("I31").Value = 2118.5667
$ws.Range("J31").Value = 5979.3
$ws.Range("K31").Value = 2118.5667
$ws.Range("L31").Value = 5979.3
$ws.Range("M31").Value = -1823.5667
$ws.Range("N31").Value = -6569.3

$ws.Range("H34").Value = 3083.75
$ws.Range("I34").Value = 2118.5667
$ws.Range("J34").Value = 5979.3
$ws.Range("K34").Value = 2118.5667
$ws.Range("L34").Value = 5979.3
$ws.Range("M34").Value = -1916.5667
$ws.Range("N34").Value = -6383.3

$ws.Range("H58").Value = 2677
$ws.Range("I58").Value = 2895.15
$ws.Range("J58").Value = 1949.8334
$ws.Range("K58").Value = 2895.15
$ws.Range("L58").Value = 1949.8334
$ws.Range("M58").Value = -2692.15
$ws.Range("N58").Value = -2355.8334

$ws.Range("H107").Value = 731.25
$ws.Range("J107").Value = 749.8889
$ws.Range("L107").Value = 749.8889
$ws.Range("N107").Value = -4589.8889

$ws.Range("H119").Value = 75000
$ws.Range("J119").Value = 75000
$ws.Range("L119").Value = 75000
$ws.Range("N119").Value = -84676

$ws.Range("H122").Value = 3598.6086
$ws.Range("I122").Value = 3530.1177
$ws.Range("K122").Value = 10590.3531
$ws.Range("M122").Value = -8140.3531

$ws.Range("H132").Value = 1363.7368
$ws.Range("I132").Value = 1392.1538
$ws.Range("K132").Value = 4176.4614
$ws.Range("M132").Value = -1646.4614

$ws.Range("H134").Value = 12684.255
$ws.Range("I134").Value = 12789.519
$ws.Range("K134").Value = 38368.557
$ws.Range("M134").Value = -35833.557

$ws.Range("H136").Value = 2677
$ws.Range("I136").Value = 2895.15
$ws.Range("J136").Value = 1949.8334
$ws.Range("K136").Value = 8685.450000000001
$ws.Range("L136").Value = 5849.5002
$ws.Range("M136").Value = -6135.450000000001
$ws.Range("N136").Value = -10949.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1598.0526
$ws.Range("I5").Value = 1218.4
$ws.Range("J5").Value = 2019.8889
$ws.Range("K5").Value = 3655.2
$ws.Range("L5").Value = 6059.6667
$ws.Range("M5").Value = -3543.2
$ws.Range("N5").Value = -6283.6667

$ws.Range("H80").Value = 1199.6666
$ws.Range("J80").Value = 1299
$ws.Range("L80").Value = 3897
$ws.Range("N80").Value = -5769

$ws.Range("H83").Value = 1199.6666
$ws.Range("J83").Value = 1299
$ws.Range("L83").Value = 11691
$ws.Range("N83").Value = -21051

$ws.Range("H107").Value = 1159.8636
$ws.Range("I107").Value = 1421.125
$ws.Range("J107").Value = 1010.5714
$ws.Range("K107").Value = 4263.375
$ws.Range("L107").Value = 3031.7142
$ws.Range("M107").Value = -2343.375
$ws.Range("N107").Value = -6871.7142

$ws.Range("H121").Value = 7343.5
$ws.Range("I121").Value = 1030
$ws.Range("J121").Value = 8606.200000000001
$ws.Range("K121").Value = 3090
$ws.Range("L121").Value = 25818.6
$ws.Range("M121").Value = -1780
$ws.Range("N121").Value = -28438.6

$ws.Range("H131").Value = 3340011.8
$ws.Range("J131").Value = 3515670.2
$ws.Range("L131").Value = 10547010.6
$ws.Range("N131").Value = -10557090.6

$ws.Range("H135").Value = 1598.0526
$ws.Range("I135").Value = 1218.4
$ws.Range("J135").Value = 2019.8889
$ws.Range("K135").Value = 10965.6
$ws.Range("L135").Value = 18179.0001
$ws.Range("M135").Value = -8430.6
$ws.Range("N135").Value = -23249.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2320
$ws.Range("J43").Value = 2500
$ws.Range("L43").Value = 2500
$ws.Range("N43").Value = -2802

$ws.Range("H51").Value = 78999.664
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 78999.664
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 78999.664
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -80017.664

$ws.Range("H70").Value = 41671356
$ws.Range("I70").Value = 83336460
$ws.Range("J70").Value = 6249.75
$ws.Range("K70").Value = 83336460
$ws.Range("L70").Value = 6249.75
$ws.Range("M70").Value = -83336190
$ws.Range("N70").Value = -6789.75

$ws.Range("H73").Value = 41671356
$ws.Range("I73").Value = 83336460
$ws.Range("J73").Value = 6249.75
$ws.Range("K73").Value = 83336460
$ws.Range("L73").Value = 6249.75
$ws.Range("M73").Value = -83335524
$ws.Range("N73").Value = -8121.75

$ws.Range("H80").Value = 5741.4614
$ws.Range("I80").Value = 4138.8335
$ws.Range("J80").Value = 7115.143
$ws.Range("K80").Value = 4138.8335
$ws.Range("L80").Value = 7115.143
$ws.Range("M80").Value = -3140.8335
$ws.Range("N80").Value = -9111.143

$ws.Range("H83").Value = 5741.4614
$ws.Range("I83").Value = 4138.8335
$ws.Range("J83").Value = 7115.143
$ws.Range("K83").Value = 20694.1675
$ws.Range("L83").Value = 35575.715
$ws.Range("M83").Value = -15702.1675
$ws.Range("N83").Value = -45559.715

$ws.Range("H99").Value = 4987.2
$ws.Range("I99").Value = 4987.2
$ws.Range("K99").Value = 4987.2
$ws.Range("M99").Value = -2741.2

$ws.Range("H102").Value = 1802.1936
$ws.Range("I102").Value = 1109.375
$ws.Range("K102").Value = 1109.375
$ws.Range("M102").Value = 512.625

$ws.Range("H107").Value = 1375.1818
$ws.Range("I107").Value = 436.1
$ws.Range("J107").Value = 2157.75
$ws.Range("K107").Value = 436.1
$ws.Range("L107").Value = 2157.75
$ws.Range("M107").Value = 1483.9
$ws.Range("N107").Value = -5997.75

$ws.Range("H132").Value = 1575.24
$ws.Range("I132").Value = 1542.0435
$ws.Range("K132").Value = 4626.1305
$ws.Range("M132").Value = -2096.1305

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1543.2142
$ws.Range("I16").Value = 1543.2142
$ws.Range("K16").Value = 1543.2142
$ws.Range("M16").Value = -1373.2142

$ws.Range("H93").Value = 4372.6
$ws.Range("I93").Value = 1665.75
$ws.Range("J93").Value = 6177.1665
$ws.Range("K93").Value = 1665.75
$ws.Range("L93").Value = 6177.1665
$ws.Range("M93").Value = -417.75
$ws.Range("N93").Value = -8673.166499999999

$ws.Range("H100").Value = 5211.05
$ws.Range("I100").Value = 3035.9
$ws.Range("J100").Value = 7386.2
$ws.Range("K100").Value = 3035.9
$ws.Range("L100").Value = 7386.2
$ws.Range("M100").Value = -2494.9
$ws.Range("N100").Value = -8468.200000000001

$ws.Range("H111").Value = 134000
$ws.Range("J111").Value = 134000
$ws.Range("L111").Value = 134000
$ws.Range("N111").Value = -142180

$ws.Range("H122").Value = 7977.4
$ws.Range("I122").Value = 8630.666999999999
$ws.Range("J122").Value = 6997.5
$ws.Range("K122").Value = 25892.001
$ws.Range("L122").Value = 20992.5
$ws.Range("M122").Value = -23442.001
$ws.Range("N122").Value = -25892.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 336.5862
$ws.Range("I107").Value = 317.69565
$ws.Range("K107").Value = 953.08695
$ws.Range("M107").Value = 966.91305

$ws.Range("H122").Value = 6506.231
$ws.Range("I122").Value = 1720.75
$ws.Range("K122").Value = 5162.25
$ws.Range("M122").Value = -2712.25

$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280
